$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above the header row (old row 2), pushing the header and
# all data rows down by one. The (now-empty) old row 1 / new row 2 stays
# blank, leaving room for a new title row at the very top.
$ws.Rows("2:2").Insert()

# Title cell "BOM" in C1, bold.
$ws.Range("C1").Value = "BOM"
$ws.Range("C1").Font.Bold = $true

# The row-insert does not repoint the existing hyperlinks, so rebuild them
# against their new (shifted down by one) cells, in the same order so the
# relationship ids line up the same way.
$links = @(
    @{ Cell = "E4"; Url = "https://github.com/jjyothilinga/PCB_AndonXbeeTerminal/" },
    @{ Cell = "E5"; Url = "https://github.com/jjyothilinga/Datasheets/" },
    @{ Cell = "E9"; Url = "https://github.com/jjyothilinga/Datasheets/" },
    @{ Cell = "E6"; Url = "https://github.com/jjyothilinga/Datasheets/" }
)
$ws.Hyperlinks.Delete()
foreach ($link in $links) {
    $ws.Hyperlinks.Add($ws.Range($link.Cell), $link.Url)
}

# Selection moves to B11 and the sheet now gets an explicit (portrait)
# page setup.
$ws.Range("B11").Select()
$ws.PageSetup.Orientation = 1
